# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> clrScheme "Office"  (used only by the Notes Master)
#   ppt/theme/theme2.xml -> clrScheme "Integral" (the deck's real/active theme,
#                            wired to Presentation + SlideMaster)
# The authored edit swaps the two themes' content wholesale, so the deck's
# visible/active theme (theme2.xml) ends up with the plain "Office" color
# palette instead of "Integral". Recolor it here via the live
# ThemeColorScheme, which is the palette actually used throughout the slides.

function Hex2Rgb($hex) {
    # PowerPoint's .RGB is a little-endian 0xBBGGRR packed value, build it
    # from a normal "RRGGBB" hex string so the mapping stays legible.
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Item(1).RGB  = Hex2Rgb "000000"  # dk1
$cs.Item(2).RGB  = Hex2Rgb "FFFFFF"  # lt1
$cs.Item(3).RGB  = Hex2Rgb "44546A"  # dk2
$cs.Item(4).RGB  = Hex2Rgb "E7E6E6"  # lt2
$cs.Item(5).RGB  = Hex2Rgb "5B9BD5"  # accent1
$cs.Item(6).RGB  = Hex2Rgb "ED7D31"  # accent2
$cs.Item(7).RGB  = Hex2Rgb "A5A5A5"  # accent3
$cs.Item(8).RGB  = Hex2Rgb "FFC000"  # accent4
$cs.Item(9).RGB  = Hex2Rgb "4472C4"  # accent5
$cs.Item(10).RGB = Hex2Rgb "70AD47"  # accent6
$cs.Item(11).RGB = Hex2Rgb "0563C1"  # hlink
$cs.Item(12).RGB = Hex2Rgb "954F72"  # folHlink
